# RNAscope: fix redundancy/typos in some reference expression matrices
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "D2_A/B/C" should be "D2_A/B/D" (duplicate of D2_C row label)
$ws.Range("A5").Value = "D2_A/B/D"

# Update the active selection/cursor position left in the sheet
$ws.Range("A6").Select() | Out-Null
